$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents (value + number format) of columns B and C for every
# data row (header row 1 through last data row 21): the "cp_site" (postal
# code) column and the "Nom" (name) column traded places.
for ($row = 1; $row -le 21; $row++) {
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)

    $bVal = $bCell.Value()
    $bFmt = $bCell.NumberFormat()
    $cVal = $cCell.Value()
    $cFmt = $cCell.NumberFormat()

    # Apply the number format before the value so numeric values that move
    # into/out of a text-formatted ("@") cell are (re)stored with the right
    # underlying cell type instead of being coerced to text.
    $bCell.NumberFormat = $cFmt
    $bCell.Value = $cVal

    $cCell.NumberFormat = $bFmt
    $cCell.Value = $bVal
}

# Column B and C widths are swapped along with their content.
$bWidth = $ws.Columns("B").Width()
$cWidth = $ws.Columns("C").Width()
$ws.Columns("B").Width = $cWidth
$ws.Columns("C").Width = $bWidth

# Update the selected cell recorded in the sheet view.
$ws.Range("D23").Select()
